# Lab01 Review Report — rebuild the review table with new "before"/"after" sections.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up cell styles by copying formatting from existing same-styled
#        cells (Range.Copy preserves both value & format; values are
#        overwritten afterwards). Do this BEFORE the old filler rows
#        (18-23) are removed, since some of the style "donors" live there.

# Row 10 needs the thin-top-border style that currently lives on row 19.
$ws.Range("B19:D19").Copy($ws.Range("B10:D10"))
# Row 11 needs the "applied empty border" style that currently lives on row 20.
$ws.Range("B20:D20").Copy($ws.Range("B11:D11"))
# Row 12 needs the plain center/wrap style that currently lives on row 21.
$ws.Range("B21:D21").Copy($ws.Range("B12:D12"))
# E/F helper columns on rows 7-12 need the plain center/wrap style (E13/F13 already carry it).
$ws.Range("E13:F13").Copy($ws.Range("E7:F7"))
$ws.Range("E13:F13").Copy($ws.Range("E8:F8"))
$ws.Range("E13:F13").Copy($ws.Range("E9:F9"))
$ws.Range("E13:F13").Copy($ws.Range("E10:F10"))
$ws.Range("E13:F13").Copy($ws.Range("E11:F11"))
$ws.Range("E13:F13").Copy($ws.Range("E12:F12"))
# Row 13 becomes the second "Calculation" header row — copy style from row 5.
$ws.Range("B5:D5").Copy($ws.Range("B13:D13"))

# --- 2. Remove the old filler rows (18-23) entirely.
$ws.Range("A18:F23").EntireRow.Delete()

# --- 3. Write the new table contents.

# "before" section
$ws.Range("A5").Value = "before"
$ws.Range("B5").Value = "Document"
$ws.Range("C5").Value = "Numar"
$ws.Range("D5").Value = "Observatii"

$ws.Range("B6").Value = "Requirements"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "no requirements at all"

$ws.Range("B7").Value = "Architecture"
$ws.Range("C7").Value = "A3"
$ws.Range("D7").Value = "not all requirements completed"

$ws.Range("B8").Value = "Code"
$ws.Range("C8").Value = "C7"
$ws.Range("D8").Value = "plaintext used, not localization"

$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "C8"
$ws.Range("D9").Value = "no input validation"

# rows 10-12 stay blank (just formatted filler rows)
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""

# "after" section
$ws.Range("A13").Value = "after"
$ws.Range("B13").Value = "Document"
$ws.Range("C13").Value = "Numar"
$ws.Range("D13").Value = "Observatii"

$ws.Range("B14").Value = "Requirements"
$ws.Range("C14").Value = "R1"
$ws.Range("D14").Value = "not all requirements implemented"

$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "R6"
$ws.Range("D15").Value = "more detail on user experience"

$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "R7"
$ws.Range("D16").Value = "requirements should be better described"

$ws.Range("B17").Value = "Architecture"
$ws.Range("C17").Value = "A10"
$ws.Range("D17").Value = "custom documents not covered by requirements"

# Rows 16 & 17 got taller (word-wrapped, longer observations) in the final file.
$ws.Range("16:16").RowHeight = 30
$ws.Range("17:17").RowHeight = 30

# --- 4. View/selection bookkeeping to match the saved workbook state.
$ws.Range("H17").Select()
$excel.ActiveWindow.WindowState = -4143
